$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.020.50"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "3.194.46"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.32"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +2.24%  "
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "3.745.50"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "60.046.29"
$ws.Range("E16").Value = "  +2.20%  "
$ws.Range("D17").Value = "3.201.44"
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.521"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.35%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "0.0₃0879"
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "2.804.91"
$ws.Range("E37").Value = "  +5.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0703"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0297"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.716"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  +3.27%  "
$ws.Range("D46").Value = "3.235.86"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.987"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.798"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.01%  "
